# Fruta / hortaliza, semanal
# This edit re-shuffles the per-row data (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg) across the
# existing data rows (2-28) of the active sheet. All other columns stay
# untouched because their values are identical on every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together as a row "record".
$cols = @("D", "J", "K", "L", "M", "P")

# Maps destination row -> source row (i.e. destination row ends up with
# the values that used to live in the source row).
$rowMap = @{
    2  = 8
    3  = 19
    4  = 10
    5  = 9
    6  = 28
    7  = 26
    8  = 18
    9  = 14
    10 = 15
    11 = 16
    12 = 21
    13 = 23
    14 = 20
    15 = 24
    16 = 2
    17 = 13
    18 = 3
    19 = 6
    20 = 7
    21 = 5
    22 = 17
    23 = 25
    24 = 12
    25 = 22
    26 = 27
    27 = 4
    28 = 11
}

# Snapshot all current values first so that writes don't clobber values
# that still need to be read for a later destination row. Value2 is used
# instead of Value because these are all plain numbers (including the
# date column, which is stored/edited as its underlying serial number).
$original = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 28; $r++) {
        $original["$col$r"] = $ws.Range("$col$r").Value2
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $original["$col$srcRow"]
    }
}
